$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Insert a new column before column H (8th column), shifting CLASSIFICATION
# and everything to its right one column over, to make room for a new
# ORGANIZATION column.
$ws.Columns.Item(8).Insert()

# Populate the new ORGANIZATION header/value column.
$ws.Cells.Item(1, 8).Value = "ORGANIZATION"
$ws.Cells.Item(2, 8).Value = "d9c76d52-03d3-4480-8c2c-b66e6d9c57f2"

# Match the width of the preceding column (G) for the newly inserted column.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Update the selected cell to reflect where the cursor ended up after edits.
$ws.Range("H2").Select()
